# Renamed sheets (and a couple columns for consistency) -- this workbook only
# touches the single worksheet's tab name: "Tabelle1" -> "Price".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Price"
